# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" rows (16-19) are reordered from descending
# (1804,1803,1802,1801) to ascending (1801,1802,1803,1804) order, while the
# "Valor Mora" figures in column F travel together with their period so the
# period -> value relationship is preserved (1801 -> 28124, 1802/1803/1804 ->
# 31249 each).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Periodo Mora (column E) - reorder to ascending period order
$ws.Range("E16").Value = "1801"
$ws.Range("E17").Value = "1802"
$ws.Range("E18").Value = "1803"
$ws.Range("E19").Value = "1804"

# Valor Mora (column F) - keep each value with its period
$ws.Range("F16").Value = 28124
$ws.Range("F17").Value = 31249
$ws.Range("F18").Value = 31249
$ws.Range("F19").Value = 31249
